$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.254.71'
$ws.Range('E2').Value = '  +5.17%  '
$ws.Range('D3').Value = '2.710.08'
$ws.Range('E3').Value = '  +4.31%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = '586.38'
$ws.Range('E5').Value = '  +0.68%  '
$ws.Range('D6').Value = '149.59'
$ws.Range('E6').Value = '  +4.72%  '
$ws.Range('D7').Value = '0.996'
$ws.Range('E7').Value = '  -0.33%  '
$ws.Range('E8').Value = '  +1.69%  '
$ws.Range('D9').Value = '2.746.87'
$ws.Range('E9').Value = '  +5.60%  '
$ws.Range('E10').Value = '  +3.16%  '
$ws.Range('E11').Value = '  +7.63%  '
$ws.Range('E12').Value = '  +4.50%  '
$ws.Range('E13').Value = '  +1.96%  '
$ws.Range('D14').Value = '3.195.80'
$ws.Range('E14').Value = '  +4.47%  '
$ws.Range('D15').Value = '26.63'
$ws.Range('E15').Value = '  +8.96%  '
$ws.Range('D16').Value = '63.141.67'
$ws.Range('E16').Value = '  +4.99%  '
$ws.Range('E17').Value = '  +7.45%  '
$ws.Range('D18').Value = '2.728.90'
$ws.Range('E18').Value = '  +4.84%  '
$ws.Range('D19').Value = '11.97'
$ws.Range('E19').Value = '  +5.67%  '
$ws.Range('D20').Value = '4.88'
$ws.Range('E20').Value = '  +5.82%  '
$ws.Range('D21').Value = '363.33'
$ws.Range('E21').Value = '  +5.19%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '7.00'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.49%  '
$ws.Range('D23').Value = '0.998'
$ws.Range('E23').Value = '  -0.20%  '
$ws.Range('D24').Value = '0.533'
$ws.Range('E24').Value = '  +0.00%  '
$ws.Range('D25').Value = '65.54'
$ws.Range('E25').Value = '  +3.02%  '
$ws.Range('E26').Value = '  +4.02%  '
$ws.Range('D27').Value = '8.65'
$ws.Range('E27').Value = '  +8.30%  '
$ws.Range('E28').Value = '  -0.23%  '
$ws.Range('D29').Value = '0.0₃0858'
$ws.Range('E29').Value = '  +7.61%  '
$ws.Range('D30').Value = '2.03'
$ws.Range('E30').Value = '  +6.35%  '
$ws.Range('D31').Value = '7.07'
$ws.Range('E31').Value = '  +10.68%  '
$ws.Range('D32').Value = '170.39'
$ws.Range('E32').Value = '  +2.13%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.20'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +22.22%  '
$ws.Range('D34').Value = '0.996'
$ws.Range('E34').Value = '  -0.21%  '
$ws.Range('E36').Value = '  +12.16%  '
$ws.Range('E37').Value = '  +8.37%  '
$ws.Range('E38').Value = '  +10.79%  '
$ws.Range('E39').Value = '  +19.10%  '
$ws.Range('D40').Value = '350.39'
$ws.Range('E40').Value = '  +11.88%  '
$ws.Range('E41').Value = '  +10.20%  '
$ws.Range('D42').Value = '39.29'
$ws.Range('E42').Value = '  +2.99%  '
$ws.Range('D43').Value = '5.67'
$ws.Range('E43').Value = '  +14.11%  '
$ws.Range('D44').Value = '21.59'
$ws.Range('E44').Value = '  +8.64%  '
$ws.Range('B45').Value = 'InjectiveProtocol'
$ws.Range('C45').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D45').Value = '21.72'
$ws.Range('E45').Value = '  +9.40%  '
$ws.Range('B46').Value = 'Hedera'
$ws.Range('C46').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D46').Value = '0.0591'
$ws.Range('E46').Value = '  +7.91%  '
$ws.Range('D47').Value = '139.33'
$ws.Range('E47').Value = '  +2.72%  '
$ws.Range('E48').Value = '  +7.37%  '
$ws.Range('D49').Value = '0.642'
$ws.Range('E49').Value = '  +6.13%  '
$ws.Range('E50').Value = '  +1.76%  '
$ws.Range('B51').Value = 'FirstDigitalUSD'
$ws.Range('C51').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D51').Value = '0.995'
$ws.Range('E51').Value = '  -0.37%  '
